# Update workbook to reflect the latest scrape:
#  1) Rows 58 and 59 had their match data (columns F:V) swapped
#     (Mornar Bar vs Mladost DG <-> Petrovac vs Decic).
#  2) A new match row (75) was appended: Sutjeska vs Petrovac.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap F:V content between row 58 and row 59 ---------------------

$row58 = $ws.Range("F58:V58").Value2
$row59 = $ws.Range("F59:V59").Value2

$ws.Range("F58:V58").Value2 = $row59
$ws.Range("F59:V59").Value2 = $row58

# --- 2) Append new row 75 ----------------------------------------------

$ws.Range("A75").Value2 = 74
$ws.Range("B75").Value = "montenegro"
$ws.Range("C75").Value = "prva-crnogorska-liga"
$ws.Range("D75").Value = "2023-2024"
$ws.Range("E75").Value2 = 45235.64583333334
$ws.Range("F75").Value = "Sutjeska"
$ws.Range("G75").Value2 = 3
$ws.Range("H75").Value = "Petrovac"
$ws.Range("I75").Value2 = 2
$ws.Range("J75").Value2 = 1.45
$ws.Range("K75").Value = "04/11/2023 03:43"
$ws.Range("L75").Value2 = 1.57
$ws.Range("M75").Value = "04/11/2023 22:30"
$ws.Range("N75").Value2 = 3.92
$ws.Range("O75").Value = "04/11/2023 03:43"
$ws.Range("P75").Value2 = 3.65
$ws.Range("Q75").Value = "04/11/2023 22:30"
$ws.Range("R75").Value2 = 5.87
$ws.Range("S75").Value = "04/11/2023 03:43"
$ws.Range("T75").Value2 = 5.93
$ws.Range("U75").Value = "04/11/2023 22:30"
$ws.Range("V75").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/sutjeska-petrovac/AmFf2v7h/"

# Apply the same style as column A's existing data cells (bold, bordered, centered)
$ws.Range("A74").Copy()
$ws.Range("A75").PasteSpecial(-4122) # xlPasteFormats
# Apply the same style as column E's existing date cells (custom date/time number format)
$ws.Range("E74").Copy()
$ws.Range("E75").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
